$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Column D (Price) updates, forced to text to match original inlineStr formatting ---
Set-TextValue "D2" "68.866.92"
Set-TextValue "D3" "3.507.48"
Set-TextValue "D5" "576.48"
Set-TextValue "D6" "173.74"
Set-TextValue "D7" "0.615"
Set-TextValue "D8" "3.497.40"
Set-TextValue "D11" "6.63"
Set-TextValue "D12" "0.602"
Set-TextValue "D13" "47.15"
Set-TextValue "D15" "678.49"
Set-TextValue "D18" "68.952.54"
Set-TextValue "D19" "3.516.40"
Set-TextValue "D21" "17.49"
Set-TextValue "D22" "11.18"
Set-TextValue "D23" "0.903"
Set-TextValue "D24" "16.37"
Set-TextValue "D25" "97.29"
Set-TextValue "D28" "2.65"
Set-TextValue "D29" "9.36"
Set-TextValue "D30" "33.01"
Set-TextValue "D31" "8.80"
Set-TextValue "D32" "3.16"
Set-TextValue "D35" "562.29"
Set-TextValue "D36" "3.64"
Set-TextValue "D37" "10.85"
Set-TextValue "D38" "0.105"
Set-TextValue "D39" "57.18"
Set-TextValue "D40" "1.00"
Set-TextValue "D42" "0.0441"
Set-TextValue "D43" "3.456.13"
Set-TextValue "D44" "0.335"
Set-TextValue "D45" "33.40"
Set-TextValue "D46" "0.0₃0701"
Set-TextValue "D47" "2.90"
Set-TextValue "D48" "2.58"
Set-TextValue "D49" "0.133"
Set-TextValue "D50" "134.26"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -4.17%  "
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -2.56%  "
$ws.Range("E6").Value = "  -4.19%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -7.25%  "
$ws.Range("E11").Value = "  +8.19%  "
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E13").Value = "  -5.64%  "
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E18").Value = "  -4.16%  "
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("E23").Value = "  -4.33%  "
$ws.Range("E24").Value = "  -9.29%  "
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("E26").Value = "  -4.98%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -6.59%  "
$ws.Range("E29").Value = "  -8.25%  "
$ws.Range("E30").Value = "  -6.67%  "
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("E32").Value = "  -9.37%  "
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("E36").Value = "  -13.26%  "
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("E42").Value = "  -6.08%  "
$ws.Range("E43").Value = "  -7.35%  "
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("E45").Value = "  -6.60%  "
$ws.Range("E46").Value = "  -8.27%  "
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("E48").Value = "  -8.02%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("E51").Value = "  -0.62%  "

# --- Row 16/17 swap: Polkadot <-> WrappedliquidstakedEther2.0 ---
Set-TextValue "B16" "WrappedliquidstakedEther2.0"
Set-TextValue "C16" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D16" "4.078.35"
$ws.Range("E16").Value = "  -3.67%  "

Set-TextValue "B17" "Polkadot"
Set-TextValue "C17" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D17" "8.88"
$ws.Range("E17").Value = "  -2.12%  "

# --- Row 33/34 swap: NEARProtocol <-> Mantle ---
Set-TextValue "B33" "Mantle"
Set-TextValue "C33" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D33" "1.35"
$ws.Range("E33").Value = "  -6.13%  "

Set-TextValue "B34" "NEARProtocol"
Set-TextValue "C34" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D34" "7.25"
$ws.Range("E34").Value = "  -1.45%  "
